$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the projected date for "Proyecto A" (row 2, column C)
$ws.Range("C2").Value = "2025-12-05"

# Clear the base/projected dates for "Proyecto B" (row 3), keep formatting
$ws.Range("B3:C3").ClearContents()

# Add a new row to the table for "Proyecto C"
$tbl = $ws.ListObjects.Item(1)
$tbl.ListRows.Add() | Out-Null
$ws.Range("A4").Value = "Proyecto C"

# New row inherits plain fill formatting (no date number format) like column A
$ws.Range("A2").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122) # xlPasteFormats

# Move the active selection
$ws.Range("I15").Select() | Out-Null
